$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SSH column (F) switches to "NO" for rows 6 and 7
$ws.Range("F6").Value = "NO"
$ws.Range("F7").Value = "NO"

# CPU-UTILS column (G) switches from "1" to "0" (as text) for rows 6, 7 and 8.
# A plain Value = "0" gets auto-coerced to a number by Excel, so we force
# text entry with a leading quote-prefix, then restore the original cell
# formatting (copied from an untouched sibling cell in the same column)
# so the style index is unaffected by the quote-prefix flag.
$ws.Range("G6").Value = "'0"
$ws.Range("G9").Copy()
$ws.Range("G6").PasteSpecial(-4122)

$ws.Range("G7").Value = "'0"
$ws.Range("G9").Copy()
$ws.Range("G7").PasteSpecial(-4122)

$ws.Range("G8").Value = "'0"
$ws.Range("G9").Copy()
$ws.Range("G8").PasteSpecial(-4122)

$excel.CutCopyMode = $false
